$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- becomes old Row4 data (A,B,E,F,G,H) with Q/R updated
$ws.Range("A2").Value = 111363016
$ws.Range("B2").Value = 77515
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 593439.5664869671
$ws.Range("R2").Value = 6986881.627536911

# Row 3 <- becomes old Row6 data (A,B,E,F,G,H) with Q/R updated
$ws.Range("A3").Value = 111363018
$ws.Range("B3").Value = 89405
$ws.Range("E3").Value = 1202
$ws.Range("F3").Value = "Ullticka"
$ws.Range("G3").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H3").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q3").Value = 593479.5069047968
$ws.Range("R3").Value = 6986870.044355935

# Row 4 <- becomes old Row3 data (A,B,E,F,G,H); Q/R unchanged
$ws.Range("A4").Value = 111363019
$ws.Range("B4").Value = 77268
$ws.Range("E4").Value = 228912
$ws.Range("F4").Value = "Mörk kolflarnlav"
$ws.Range("G4").Value = "Carbonicola myrmecina"
$ws.Range("H4").Value = "(Ach.) Bendiksby & Timdal"

# Row 6 <- becomes old Row2 data (A,B,E,F,G,H) with Q/R updated
$ws.Range("A6").Value = 111363017
$ws.Range("B6").Value = 78578
$ws.Range("E6").Value = 6458
$ws.Range("F6").Value = "Lunglav"
$ws.Range("G6").Value = "Lobaria pulmonaria"
$ws.Range("H6").Value = "(L.) Hoffm."
$ws.Range("Q6").Value = 593472.3298762256
$ws.Range("R6").Value = 6986898.025413335
